# Scheduled runner update: refresh market-board derived profit figures
# across the per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 408.5
$ws.Range("I11").Value = 408.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 408.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -268.5
$ws.Range("H12").Value = 1375
$ws.Range("I12").Value = 1375
$ws.Range("J12").Value = 1375
$ws.Range("K12").Value = 1375
$ws.Range("L12").Value = 1375
$ws.Range("M12").Value = -1205
$ws.Range("N12").Value = -1715
$ws.Range("H53").Value = 404.15
$ws.Range("I53").Value = 486.5
$ws.Range("J53").Value = 212
$ws.Range("K53").Value = 486.5
$ws.Range("L53").Value = 212
$ws.Range("M53").Value = 150.5
$ws.Range("N53").Value = -1486
$ws.Range("H69").Value = 6197.1665
$ws.Range("I69").Value = 3950
$ws.Range("J69").Value = 6478.0625
$ws.Range("K69").Value = 11850
$ws.Range("L69").Value = 19434.1875
$ws.Range("M69").Value = -10976
$ws.Range("N69").Value = -21182.1875
$ws.Range("H72").Value = 6197.1665
$ws.Range("I72").Value = 3950
$ws.Range("J72").Value = 6478.0625
$ws.Range("K72").Value = 35550
$ws.Range("L72").Value = 58302.5625
$ws.Range("M72").Value = -31182
$ws.Range("N72").Value = -67038.5625
$ws.Range("H107").Value = 525
$ws.Range("I107").Value = 523
$ws.Range("J107").Value = 555
$ws.Range("K107").Value = 523
$ws.Range("L107").Value = 555
$ws.Range("M107").Value = 1397
$ws.Range("N107").Value = -4395
$ws.Range("H112").Value = 945.4103
$ws.Range("I112").Value = 1800
$ws.Range("J112").Value = 922.9211
$ws.Range("K112").Value = 5400
$ws.Range("L112").Value = 2768.7633
$ws.Range("M112").Value = -4292
$ws.Range("N112").Value = -4984.763300000001
$ws.Range("H113").Value = 4534.8
$ws.Range("I113").Value = 3558
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 3558
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -304
$ws.Range("N113").Value = -12508
$ws.Range("H118").Value = 716.2
$ws.Range("I118").Value = 716.2
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2148.6
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -491.6000000000004
$ws.Range("H132").Value = 2671.535
$ws.Range("I132").Value = 861.9729599999999
$ws.Range("J132").Value = 13830.5
$ws.Range("K132").Value = 2585.91888
$ws.Range("L132").Value = 41491.5
$ws.Range("M132").Value = -55.91887999999972
$ws.Range("H137").Value = 1997.3715
$ws.Range("I137").Value = 1095.5416
$ws.Range("J137").Value = 3965
$ws.Range("K137").Value = 3286.6248
$ws.Range("L137").Value = 11895
$ws.Range("M137").Value = -736.6248000000001
$ws.Range("H138").Value = 3192.1343
$ws.Range("I138").Value = 1583.2858
$ws.Range("J138").Value = 3926.6086
$ws.Range("K138").Value = 4749.857400000001
$ws.Range("L138").Value = 11779.8258
$ws.Range("M138").Value = 390.1425999999992
$ws.Range("N138").Value = -22059.8258
$ws.Range("H141").Value = 1147.5
$ws.Range("I141").Value = 1147.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3442.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1737.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1761.707
$ws.Range("I32").Value = 1070.7831
$ws.Range("J32").Value = 5345.875
$ws.Range("K32").Value = 1070.7831
$ws.Range("L32").Value = 5345.875
$ws.Range("M32").Value = -783.7831000000001
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H61").Value = 1791.3671
$ws.Range("I61").Value = 1484.4546
$ws.Range("J61").Value = 3349.5386
$ws.Range("K61").Value = 1484.4546
$ws.Range("L61").Value = 3349.5386
$ws.Range("M61").Value = -1272.4546
$ws.Range("H74").Value = 2456.5862
$ws.Range("I74").Value = 1464.55
$ws.Range("J74").Value = 4661.1113
$ws.Range("K74").Value = 1464.55
$ws.Range("L74").Value = 4661.1113
$ws.Range("M74").Value = -590.55
$ws.Range("H77").Value = 2456.5862
$ws.Range("I77").Value = 1464.55
$ws.Range("J77").Value = 4661.1113
$ws.Range("K77").Value = 7322.75
$ws.Range("L77").Value = 23305.5565
$ws.Range("M77").Value = -2954.75
$ws.Range("H132").Value = 2424.4211
$ws.Range("I132").Value = 2282.9707
$ws.Range("J132").Value = 3626.75
$ws.Range("K132").Value = 6848.9121
$ws.Range("L132").Value = 10880.25
$ws.Range("M132").Value = -4318.9121
$ws.Range("H136").Value = 1791.3671
$ws.Range("I136").Value = 1484.4546
$ws.Range("J136").Value = 3349.5386
$ws.Range("K136").Value = 4453.3638
$ws.Range("L136").Value = 10048.6158
$ws.Range("M136").Value = -1903.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1920.8572
$ws.Range("I80").Value = 528
$ws.Range("J80").Value = 2478
$ws.Range("K80").Value = 528
$ws.Range("L80").Value = 2478
$ws.Range("M80").Value = 470
$ws.Range("N80").Value = -4474
$ws.Range("H83").Value = 1920.8572
$ws.Range("I83").Value = 528
$ws.Range("J83").Value = 2478
$ws.Range("K83").Value = 2640
$ws.Range("L83").Value = 12390
$ws.Range("M83").Value = 2352
$ws.Range("N83").Value = -22374
$ws.Range("H86").Value = 1774.8182
$ws.Range("I86").Value = 1734.8334
$ws.Range("J86").Value = 1822.8
$ws.Range("K86").Value = 1734.8334
$ws.Range("L86").Value = 1822.8
$ws.Range("M86").Value = -611.8334
$ws.Range("N86").Value = -4068.8
$ws.Range("H89").Value = 1774.8182
$ws.Range("I89").Value = 1734.8334
$ws.Range("J89").Value = 1822.8
$ws.Range("K89").Value = 8674.166999999999
$ws.Range("L89").Value = 9114
$ws.Range("M89").Value = -3058.166999999999
$ws.Range("N89").Value = -20346
$ws.Range("H99").Value = 34616.5
$ws.Range("I99").Value = 36019.25
$ws.Range("J99").Value = 29005.5
$ws.Range("K99").Value = 36019.25
$ws.Range("L99").Value = 29005.5
$ws.Range("M99").Value = -34521.25
$ws.Range("N99").Value = -32001.5
$ws.Range("H134").Value = 3210.8823
$ws.Range("I134").Value = 2670.4783
$ws.Range("J134").Value = 8182.6
$ws.Range("K134").Value = 8011.4349
$ws.Range("L134").Value = 24547.8
$ws.Range("M134").Value = -5476.4349
$ws.Range("N134").Value = -29617.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6443.2573
$ws.Range("I31").Value = 5601
$ws.Range("J31").Value = 8057.5835
$ws.Range("K31").Value = 5601
$ws.Range("L31").Value = 8057.5835
$ws.Range("M31").Value = -5306
$ws.Range("H34").Value = 6443.2573
$ws.Range("I34").Value = 5601
$ws.Range("J34").Value = 8057.5835
$ws.Range("K34").Value = 5601
$ws.Range("L34").Value = 8057.5835
$ws.Range("M34").Value = -5399
$ws.Range("H62").Value = 2615.4167
$ws.Range("I62").Value = 2498.5
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 2498.5
$ws.Range("L62").Value = 3200
$ws.Range("M62").Value = -1874.5
$ws.Range("N62").Value = -4448
$ws.Range("H65").Value = 2615.4167
$ws.Range("I65").Value = 2498.5
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 12492.5
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = -9372.5
$ws.Range("N65").Value = -22240
$ws.Range("H86").Value = 1999.5
$ws.Range("I86").Value = 1999.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1999.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -876.5
$ws.Range("H89").Value = 1999.5
$ws.Range("I89").Value = 1999.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9997.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4381.5
$ws.Range("H132").Value = 167482.06
$ws.Range("I132").Value = 766.5
$ws.Range("J132").Value = 2501500
$ws.Range("K132").Value = 2299.5
$ws.Range("L132").Value = 7504500
$ws.Range("M132").Value = 230.5
$ws.Range("H134").Value = 7286
$ws.Range("I134").Value = 7502.353
$ws.Range("J134").Value = 5447
$ws.Range("K134").Value = 22507.059
$ws.Range("L134").Value = 16341
$ws.Range("M134").Value = -19972.059

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 145971.86
$ws.Range("I139").Value = 202681.8
$ws.Range("J139").Value = 4197
$ws.Range("K139").Value = 608045.3999999999
$ws.Range("L139").Value = 12591
$ws.Range("M139").Value = -602905.3999999999
$ws.Range("H140").Value = 3852.4375
$ws.Range("I140").Value = 3780.1667
$ws.Range("J140").Value = 3895.8
$ws.Range("K140").Value = 11340.5001
$ws.Range("L140").Value = 11687.4
$ws.Range("M140").Value = -6160.500100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2910
$ws.Range("I126").Value = 2800
$ws.Range("J126").Value = 3350
$ws.Range("K126").Value = 8400
$ws.Range("L126").Value = 10050
$ws.Range("M126").Value = -5930
$ws.Range("H132").Value = 22232754
$ws.Range("I132").Value = 23817190
$ws.Range("J132").Value = 50640.332
$ws.Range("K132").Value = 71451570
$ws.Range("L132").Value = 151920.996
$ws.Range("M132").Value = -71449040

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4111
$ws.Range("I107").Value = 3360
$ws.Range("J107").Value = 4937.1
$ws.Range("K107").Value = 10080
$ws.Range("L107").Value = 14811.3
$ws.Range("M107").Value = -8160
$ws.Range("N107").Value = -18651.3
$ws.Range("H113").Value = 2114.8333
$ws.Range("I113").Value = 959.1875
$ws.Range("J113").Value = 4426.125
$ws.Range("K113").Value = 2877.5625
$ws.Range("L113").Value = 13278.375
$ws.Range("M113").Value = -707.5625
$ws.Range("N113").Value = -17618.375
$ws.Range("H132").Value = 2221.8276
$ws.Range("I132").Value = 1713.2
$ws.Range("J132").Value = 5400.75
$ws.Range("K132").Value = 5139.6
$ws.Range("L132").Value = 16202.25
$ws.Range("M132").Value = -2609.6
$ws.Range("H136").Value = 2770.6206
$ws.Range("I136").Value = 2290.434
$ws.Range("J136").Value = 7860.6
$ws.Range("K136").Value = 6871.302000000001
$ws.Range("L136").Value = 23581.8
$ws.Range("M136").Value = -4321.302000000001

